# Week 13 logging update
# Applies new game/play data to the YDS, OFF, DEF, ST, TURNS and PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's play-by-play yardage figures to the
# running space-separated lists stored in B2/C2/B3/C3 (OFF/DEF plays).
# ---------------------------------------------------------------------------
$yds = $wb.Worksheets("YDS")

$b2 = $yds.Range("B2").Value()
$yds.Range("B2").Value = $b2 + " 2 5 3 -1 3 2 9 9 3 -1 25 24 7 0 4 8 6 5 9 23 4 6 8 0 7 7"

$c2 = $yds.Range("C2").Value()
$yds.Range("C2").Value = $c2 + " 0 3 2 1 2 33 2 1 3 4 -2 4 58 2 10 2 3 5 3 3 10"

$b3 = $yds.Range("B3").Value()
$yds.Range("B3").Value = $b3 + " 6 12 -2 -6 25 27 9 24 9 1 9 30 6 1 5 13 13 70"

$c3 = $yds.Range("C3").Value()
$yds.Range("C3").Value = $c3 + " 9 5 1 41 7 1 21 10 10 5 17 1 4 20 0 12 2 10 1 3 4 13 25 9 4 3"

# ---------------------------------------------------------------------------
# OFF sheet: updated season totals for Home (row 2) and Road (row 3).
# ---------------------------------------------------------------------------
$off = $wb.Worksheets("OFF")

$off.Range("C2").Value = 144
$off.Range("E2").Value = 9
$off.Range("F2").Value = 55
$off.Range("G2").Value = 51
$off.Range("I2").Value = 10
$off.Range("J2").Value = 20
$off.Range("L2").Value = 199
$off.Range("M2").Value = 117
$off.Range("O2").Value = 16
$off.Range("P2").Value = 9
$off.Range("Q2").Value = 382

$off.Range("B3").Value = 10
$off.Range("C3").Value = 127
$off.Range("F3").Value = 74
$off.Range("G3").Value = 16
$off.Range("H3").Value = 22
$off.Range("I3").Value = 42
$off.Range("J3").Value = 34
$off.Range("N3").Value = 6

# ---------------------------------------------------------------------------
# DEF sheet: updated season totals for Home (row 2) and Road (row 3).
# ---------------------------------------------------------------------------
$def = $wb.Worksheets("DEF")

$def.Range("C2").Value = 145
$def.Range("F2").Value = 46
$def.Range("G2").Value = 35
$def.Range("L2").Value = 206
$def.Range("M2").Value = 142
$def.Range("O2").Value = 13
$def.Range("Q2").Value = 335

$def.Range("B3").Value = 9
$def.Range("C3").Value = 122
$def.Range("D3").Value = 2
$def.Range("E3").Value = 26
$def.Range("F3").Value = 81
$def.Range("G3").Value = 19
$def.Range("I3").Value = 51
$def.Range("J3").Value = 41
$def.Range("N3").Value = 9

# ---------------------------------------------------------------------------
# ST sheet: updated special teams totals (row 2) and field goal attempts
# (row 3), plus this week's per-game figures appended to the running lists
# in B4/B5 (RA/RM counts) and D3/D4/D5 (TB/RA/RM splits).
# ---------------------------------------------------------------------------
$st = $wb.Worksheets("ST")

$st.Range("B2").Value = 55
$st.Range("D2").Value = 54
$st.Range("F2").Value = 77
$st.Range("G2").Value = 74
$st.Range("L2").Value = 17
$st.Range("M2").Value = 9
$st.Range("N2").Value = 16
$st.Range("B3").Value = 37

$stB4 = $st.Range("B4").Value()
$st.Range("B4").Value = $stB4 + " 68"

$stB5 = $st.Range("B5").Value()
$st.Range("B5").Value = $stB5 + " 22"

$stD3 = $st.Range("D3").Value()
$st.Range("D3").Value = $stD3 + " 47 59 52 55 33 62"

$stD4 = $st.Range("D4").Value()
$st.Range("D4").Value = $stD4 + " 7 16 9 0 0 12"

$stD5 = $st.Range("D5").Value()
$st.Range("D5").Value = $stD5 + " 8 19 0 8 0 2 0"

# ---------------------------------------------------------------------------
# TURNS sheet: updated Home turnover totals (row 2).
# ---------------------------------------------------------------------------
$turns = $wb.Worksheets("TURNS")

$turns.Range("B2").Value = 7
$turns.Range("C2").Value = 8
$turns.Range("D2").Value = 2

# ---------------------------------------------------------------------------
# PEN sheet: updated penalty counts.
# ---------------------------------------------------------------------------
$pen = $wb.Worksheets("PEN")

$pen.Range("B2").Value = 14
$pen.Range("B3").Value = 15
$pen.Range("D3").Value = 2
$pen.Range("D4").Value = 8
